# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Betarraga"
# at row 470 (pushing the existing rows 470:487 down to 471:488), then populate
# the newly inserted row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 470; Excel shifts row 470:487 -> 471:488
# and the new blank row inherits formatting (date style) from the row above, as
# a native "Insert" does.
$ws.Rows("470:470").Insert()

# Populate the newly inserted row 470 with the new weekly record.
$ws.Cells.Item(470, 1).Value = 4
$ws.Cells.Item(470, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(470, 3).Value = "Los Lagos"
$ws.Cells.Item(470, 4).Value = 45075
$ws.Cells.Item(470, 5).Value = 10
$ws.Cells.Item(470, 6).Value = 100114014
$ws.Cells.Item(470, 7).Value = "Betarraga"
$ws.Cells.Item(470, 8).Value = "Sin especificar"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 250
$ws.Cells.Item(470, 11).Value = 1100
$ws.Cells.Item(470, 12).Value = 1100
$ws.Cells.Item(470, 13).Value = 1100
$ws.Cells.Item(470, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(470, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(470, 16).Value = 220
$ws.Cells.Item(470, 17).Value = 5
$ws.Cells.Item(470, 18).Value = "Hortaliza"
